$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("50").Insert()

$ws.Range("A50").Value = 11
$ws.Range("B50").Value = "Vega Monumental Concepción"
$ws.Range("C50").Value = "Bíobío"
$ws.Range("D50").Value = 45211
$ws.Range("E50").Value = 8
$ws.Range("F50").Value = "Fruta"
$ws.Range("G50").Value = 100107
$ws.Range("H50").Value = "Otros"
$ws.Range("I50").Value = 100107002
$ws.Range("J50").Value = "Chirimoya"
$ws.Range("K50").Value = "Cultivar IV Región"
$ws.Range("L50").Value = "Primera"
$ws.Range("M50").Value = 100
$ws.Range("N50").Value = 17000
$ws.Range("O50").Value = 21000
$ws.Range("P50").Value = 19000
$ws.Range("Q50").Value = "$/bandeja 10 kilos"
$ws.Range("R50").Value = "Provincia de Limarí"
$ws.Range("S50").Value = 1900
$ws.Range("T50").Value = 10
